# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column A ("numero_imovel") before the current "uf" column.
# ---------------------------------------------------------------------------
$ws.Range("A:A").Insert()
$ws.Range("A1").Value = "numero_imovel"
# Copy the header style (bold font, border, centered) from the neighbouring
# header cell so the new header matches the rest of the header row exactly.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Insert two new columns ("modalidade_venda", "imovel") right before the
#    current "metro_quadrado_m3" column (which, after step 1, sits at H).
#    Excel extends the header formatting automatically here because both
#    neighbouring header cells (G1/J1) already carry it.
# ---------------------------------------------------------------------------
$ws.Range("H:I").Insert()
$ws.Range("H1").Value = "modalidade_venda"
$ws.Range("I1").Value = "imovel"

# ---------------------------------------------------------------------------
# 3) Append a new trailing column ("data_processamento") after "link_acesso".
# ---------------------------------------------------------------------------
$ws.Range("M1").Value = "data_processamento"
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Fill in the data for the new columns, row by row.
# ---------------------------------------------------------------------------
$numeroImovel = @{2="10137825"; 3="10137826"; 4="10137827"; 5="1444407695527"; 6="10137828"; 7="10171235"; 8="10153991"; 9="8444409794539"; 10="10137830"; 11="10153988"}
$modalidade   = @{2="Venda Online"; 3="Venda Online"; 4="Venda Online"; 5="Venda Direta Online"; 6="Venda Online"; 7="Venda Direta Online"; 8="Venda Online"; 9="Venda Direta Online"; 10="Venda Online"; 11="Venda Online"}
$imovelTipo   = @{2="Terreno"; 3="Terreno"; 4="Terreno"; 5="Terreno"; 6="Terreno"; 7="Terreno"; 8="Terreno"; 9="Casa"; 10="Terreno"; 11="Terreno"}

foreach ($r in 2..11) {
    $ws.Cells.Item($r, 1).Value = [double]$numeroImovel[$r]
    $ws.Cells.Item($r, 8).Value = $modalidade[$r]
    $ws.Cells.Item($r, 9).Value = $imovelTipo[$r]
    $ws.Cells.Item($r, 13).Value = 45673.75602836806
}

# ---------------------------------------------------------------------------
# 5) Apply the "data_processamento" timestamp number format. Set it on M2
#    first (registers both the intermediate lowercase format and the final
#    uppercase format in numFmts, mirroring the original file), then copy
#    that cell's format down to the rest of the column in one shot so every
#    data row shares the same style index instead of minting a new one per
#    write.
# ---------------------------------------------------------------------------
$ws.Range("M2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("M2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M2").Copy()
$ws.Range("M3:M11").PasteSpecial(-4122)
